$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / header text updates (October -> November)
$ws.Range("A2").Value = "Year-to-Date through November 2016 and November 2015 (Million Cubic Feet)"
$ws.Range("B5").Value = "November 2016 YTD"
$ws.Range("E5").Value = "November 2016 YTD"
$ws.Range("G5").Value = "November 2016 YTD"
$ws.Range("I5").Value = "November 2016 YTD"
$ws.Range("K5").Value = "November 2016 YTD"
$ws.Range("C5").Value = "November 2015 YTD"
$ws.Range("F5").Value = "November 2015 YTD"
$ws.Range("H5").Value = "November 2015 YTD"
$ws.Range("J5").Value = "November 2015 YTD"
$ws.Range("L5").Value = "November 2015 YTD"

# Data cell updates
# Row 6
$ws.Range("B6").Value = 11680
$ws.Range("C6").Value = 10508
$ws.Range("D6").Value = 0.11
$ws.Range("G6").Value = 10934
$ws.Range("H6").Value = 10006
$ws.Range("I6").Value = 746
$ws.Range("J6").Value = 502

# Row 7
$ws.Range("B7").Value = 493
$ws.Range("C7").Value = 445
$ws.Range("D7").Value = 0.11
$ws.Range("G7").Value = 493
$ws.Range("H7").Value = 445

# Row 8
$ws.Range("B8").Value = 867
$ws.Range("C8").Value = 773
$ws.Range("D8").Value = 0.12
$ws.Range("G8").Value = 867
$ws.Range("H8").Value = 773

# Row 9
$ws.Range("B9").Value = 4113
$ws.Range("C9").Value = 3477
$ws.Range("D9").Value = 0.18
$ws.Range("G9").Value = 4113
$ws.Range("H9").Value = 3477

# Row 10
$ws.Range("B10").Value = 1937
$ws.Range("C10").Value = 1493
$ws.Range("G10").Value = 1190
$ws.Range("H10").Value = 992
$ws.Range("I10").Value = 746
$ws.Range("J10").Value = 502

# Row 11
$ws.Range("B11").Value = 3783
$ws.Range("C11").Value = 3827
$ws.Range("D11").Value = -0.011
$ws.Range("G11").Value = 3783
$ws.Range("H11").Value = 3827

# Row 12
$ws.Range("B12").Value = 487
$ws.Range("C12").Value = 493
$ws.Range("D12").Value = -0.012
$ws.Range("G12").Value = 487
$ws.Range("H12").Value = 493

# Row 13
$ws.Range("B13").Value = 56582
$ws.Range("C13").Value = 50424
$ws.Range("D13").Value = 0.12
$ws.Range("G13").Value = 53998
$ws.Range("H13").Value = 48160
$ws.Range("I13").Value = 1058
$ws.Range("J13").Value = 749
$ws.Range("K13").Value = 1526
$ws.Range("L13").Value = 1515

# Row 14
$ws.Range("B14").Value = 9575
$ws.Range("C14").Value = 8250
$ws.Range("D14").Value = 0.16
$ws.Range("G14").Value = 9176
$ws.Range("H14").Value = 7976
$ws.Range("I14").Value = 400
$ws.Range("J14").Value = 274

# Row 15
$ws.Range("B15").Value = 17187
$ws.Range("C15").Value = 14650
$ws.Range("D15").Value = 0.17
$ws.Range("G15").Value = 17187
$ws.Range("H15").Value = 14650

# Row 16
$ws.Range("B16").Value = 29819
$ws.Range("C16").Value = 27524
$ws.Range("D16").Value = 0.083
$ws.Range("G16").Value = 27635
$ws.Range("H16").Value = 25534
$ws.Range("I16").Value = 659
$ws.Range("J16").Value = 475
$ws.Range("K16").Value = 1526
$ws.Range("L16").Value = 1515

# Row 17
$ws.Range("B17").Value = 72736
$ws.Range("C17").Value = 59088
$ws.Range("E17").Value = 8610
$ws.Range("F17").Value = 6546
$ws.Range("G17").Value = 63145
$ws.Range("H17").Value = 51616
$ws.Range("I17").Value = 250
$ws.Range("J17").Value = 208
$ws.Range("K17").Value = 732
$ws.Range("L17").Value = 719

# Row 18
$ws.Range("B18").Value = 15451
$ws.Range("C18").Value = 12373
$ws.Range("D18").Value = 0.25
$ws.Range("E18").Value = 455
$ws.Range("F18").Value = 337
$ws.Range("G18").Value = 14997
$ws.Range("H18").Value = 12036

# Row 19
$ws.Range("B19").Value = 9436
$ws.Range("C19").Value = 7265
$ws.Range("E19").Value = 7932
$ws.Range("F19").Value = 6022
$ws.Range("G19").Value = 1204
$ws.Range("H19").Value = 986
$ws.Range("K19").Value = 300
$ws.Range("L19").Value = 257

# Row 20
$ws.Range("B20").Value = 22078
$ws.Range("C20").Value = 18224
$ws.Range("G20").Value = 22078
$ws.Range("H20").Value = 18224

# Row 21
$ws.Range("B21").Value = 12527
$ws.Range("C21").Value = 10359
$ws.Range("G21").Value = 12527
$ws.Range("H21").Value = 10359

# Row 22
$ws.Range("B22").Value = 13244
$ws.Range("C22").Value = 10868
$ws.Range("D22").Value = 0.22
$ws.Range("E22").Value = 223
$ws.Range("F22").Value = 187
$ws.Range("G22").Value = 12339
$ws.Range("H22").Value = 10011
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 208
$ws.Range("K22").Value = 432
$ws.Range("L22").Value = 462

# Row 23
$ws.Range("B23").Value = 12768
$ws.Range("C23").Value = 9939
$ws.Range("D23").Value = 0.28
$ws.Range("E23").Value = 3664
$ws.Range("F23").Value = 2692
$ws.Range("G23").Value = 9105
$ws.Range("H23").Value = 7247

# Row 24
$ws.Range("B24").Value = 3097
$ws.Range("C24").Value = 2460
$ws.Range("G24").Value = 3097
$ws.Range("H24").Value = 2460

# Row 25
$ws.Range("B25").Value = 1626
$ws.Range("C25").Value = 1258
$ws.Range("D25").Value = 0.29
$ws.Range("G25").Value = 1626
$ws.Range("H25").Value = 1258

# Row 26
$ws.Range("B26").Value = 4059
$ws.Range("C26").Value = 3212
$ws.Range("D26").Value = 0.26
$ws.Range("E26").Value = 932
$ws.Range("F26").Value = 668
$ws.Range("G26").Value = 3127
$ws.Range("H26").Value = 2544

# Row 27
$ws.Range("B27").Value = 2316
$ws.Range("C27").Value = 1815
$ws.Range("D27").Value = 0.28
$ws.Range("E27").Value = 1060
$ws.Range("F27").Value = 829
$ws.Range("G27").Value = 1255
$ws.Range("H27").Value = 986

# Row 28
$ws.Range("B28").Value = 1671
$ws.Range("C28").Value = 1194
$ws.Range("D28").Value = 0.4
$ws.Range("E28").Value = 1671
$ws.Range("F28").Value = 1194

# Row 31
$ws.Range("B31").Value = 51326
$ws.Range("C31").Value = 42633
$ws.Range("D31").Value = 0.2
$ws.Range("E31").Value = 5237
$ws.Range("F31").Value = 5038
$ws.Range("G31").Value = 40898
$ws.Range("H31").Value = 33714
$ws.Range("I31").Value = 2659
$ws.Range("J31").Value = 1911
$ws.Range("K31").Value = 2531
$ws.Range("L31").Value = 1970

# Row 32
$ws.Range("B32").Value = 1721
$ws.Range("C32").Value = 1439
$ws.Range("D32").Value = 0.2
$ws.Range("G32").Value = 1513
$ws.Range("H32").Value = 1295
$ws.Range("K32").Value = 208
$ws.Range("L32").Value = 144

# Row 34
$ws.Range("B34").Value = 6816
$ws.Range("C34").Value = 7985
$ws.Range("D34").Value = -0.15
$ws.Range("E34").Value = 1574
$ws.Range("F34").Value = 1906
$ws.Range("G34").Value = 5205
$ws.Range("H34").Value = 6044
$ws.Range("J34").Value = 7
$ws.Range("L34").Value = 28

# Row 35
$ws.Range("B35").Value = 4893
$ws.Range("C35").Value = 3927
$ws.Range("D35").Value = 0.25
$ws.Range("G35").Value = 4355
$ws.Range("H35").Value = 3448
$ws.Range("K35").Value = 538
$ws.Range("L35").Value = 402

# Row 36
$ws.Range("B36").Value = 2575
$ws.Range("C36").Value = 2200
$ws.Range("D36").Value = 0.17
$ws.Range("G36").Value = 1793
$ws.Range("H36").Value = 1506
$ws.Range("I36").Value = 782
$ws.Range("J36").Value = 693

# Row 37
$ws.Range("B37").Value = 12532
$ws.Range("C37").Value = 8979
$ws.Range("D37").Value = 0.4
$ws.Range("G37").Value = 10978
$ws.Range("H37").Value = 8086
$ws.Range("I37").Value = 1555
$ws.Range("J37").Value = 893

# Row 38
$ws.Range("B38").Value = 5745
$ws.Range("C38").Value = 4774
$ws.Range("D38").Value = 0.2
$ws.Range("E38").Value = 3578
$ws.Range("F38").Value = 3065
$ws.Range("G38").Value = 409
$ws.Range("H38").Value = 314
$ws.Range("K38").Value = 1758
$ws.Range("L38").Value = 1396

# Row 39
$ws.Range("B39").Value = 16888
$ws.Range("C39").Value = 13194
$ws.Range("D39").Value = 0.28
$ws.Range("F39").Value = 66
$ws.Range("G39").Value = 16490
$ws.Range("H39").Value = 12886
$ws.Range("I39").Value = 312
$ws.Range("J39").Value = 241

# Row 40
$ws.Range("B40").Value = 156
$ws.Range("C40").Value = 136
$ws.Range("D40").Value = 0.15
$ws.Range("G40").Value = 156
$ws.Range("H40").Value = 136

# Row 41
$ws.Range("B41").Value = 5740
$ws.Range("C41").Value = 4709
$ws.Range("D41").Value = 0.22
$ws.Range("E41").Value = 2265
$ws.Range("F41").Value = 1899
$ws.Range("G41").Value = 3475
$ws.Range("H41").Value = 2810

# Row 42
$ws.Range("B42").Value = 1133
$ws.Range("C42").Value = 958
$ws.Range("D42").Value = 0.18
$ws.Range("G42").Value = 1133
$ws.Range("H42").Value = 958

# Row 43
$ws.Range("B43").Value = 2442
$ws.Range("C43").Value = 2042
$ws.Range("D43").Value = 0.2
$ws.Range("E43").Value = 2265
$ws.Range("F43").Value = 1899
$ws.Range("G43").Value = 177
$ws.Range("H43").Value = 142

# Row 44
$ws.Range("B44").Value = 277
$ws.Range("C44").Value = 210
$ws.Range("D44").Value = 0.32
$ws.Range("G44").Value = 277
$ws.Range("H44").Value = 210

# Row 45
$ws.Range("B45").Value = 1889
$ws.Range("C45").Value = 1499
$ws.Range("D45").Value = 0.26
$ws.Range("G45").Value = 1889
$ws.Range("H45").Value = 1499

# Row 46
$ws.Range("B46").Value = 19742
$ws.Range("C46").Value = 15555
$ws.Range("D46").Value = 0.27
$ws.Range("G46").Value = 18945
$ws.Range("H46").Value = 14975
$ws.Range("I46").Value = 797
$ws.Range("J46").Value = 580

# Row 47
$ws.Range("B47").Value = 1858
$ws.Range("C47").Value = 1436
$ws.Range("D47").Value = 0.29
$ws.Range("G47").Value = 1858
$ws.Range("H47").Value = 1436

# Row 49
$ws.Range("B49").Value = 350
$ws.Range("C49").Value = 269
$ws.Range("D49").Value = 0.3
$ws.Range("G49").Value = 350
$ws.Range("H49").Value = 269

# Row 50
$ws.Range("B50").Value = 17535
$ws.Range("C50").Value = 13850
$ws.Range("D50").Value = 0.27
$ws.Range("G50").Value = 16738
$ws.Range("H50").Value = 13270
$ws.Range("I50").Value = 797
$ws.Range("J50").Value = 580

# Row 51
$ws.Range("B51").Value = 5839
$ws.Range("C51").Value = 5291
$ws.Range("D51").Value = 0.1
$ws.Range("E51").Value = 319
$ws.Range("F51").Value = 535
$ws.Range("G51").Value = 4922
$ws.Range("H51").Value = 4283
$ws.Range("I51").Value = 599
$ws.Range("J51").Value = 473

# Row 52
$ws.Range("B52").Value = 1080
$ws.Range("C52").Value = 1190
$ws.Range("D52").Value = -0.093
$ws.Range("G52").Value = 1080
$ws.Range("H52").Value = 909

# Row 53
$ws.Range("B53").Value = 1208
$ws.Range("C53").Value = 1160
$ws.Range("D53").Value = 0.042
$ws.Range("G53").Value = 1208
$ws.Range("H53").Value = 1160

# Row 54
$ws.Range("B54").Value = 1035
$ws.Range("C54").Value = 905
$ws.Range("D54").Value = 0.14
$ws.Range("E54").Value = 319
$ws.Range("F54").Value = 254
$ws.Range("G54").Value = 498
$ws.Range("H54").Value = 500
$ws.Range("I54").Value = 218
$ws.Range("J54").Value = 151

# Row 56
$ws.Range("B56").Value = 620
$ws.Range("C56").Value = 495
$ws.Range("D56").Value = 0.25
$ws.Range("G56").Value = 620
$ws.Range("H56").Value = 495

# Row 58
$ws.Range("B58").Value = 1897
$ws.Range("C58").Value = 1540
$ws.Range("D58").Value = 0.23
$ws.Range("G58").Value = 1517
$ws.Range("H58").Value = 1219
$ws.Range("I58").Value = 380
$ws.Range("J58").Value = 322

# Row 60
$ws.Range("B60").Value = 63621
$ws.Range("C60").Value = 57848
$ws.Range("D60").Value = 0.1
$ws.Range("E60").Value = 7890
$ws.Range("F60").Value = 6382
$ws.Range("G60").Value = 40049
$ws.Range("H60").Value = 33768
$ws.Range("I60").Value = 15683
$ws.Range("J60").Value = 17698

# Row 61
$ws.Range("B61").Value = 52231
$ws.Range("C61").Value = 48864
$ws.Range("D61").Value = 0.069
$ws.Range("E61").Value = 2627
$ws.Range("F61").Value = 2368
$ws.Range("G61").Value = 34432
$ws.Range("H61").Value = 29218
$ws.Range("I61").Value = 15172
$ws.Range("J61").Value = 17277

# Row 62
$ws.Range("B62").Value = 6252
$ws.Range("C62").Value = 5186
$ws.Range("E62").Value = 1486
$ws.Range("F62").Value = 1391
$ws.Range("G62").Value = 4256
$ws.Range("H62").Value = 3374
$ws.Range("I62").Value = 511
$ws.Range("J62").Value = 421

# Row 63
$ws.Range("B63").Value = 5138
$ws.Range("C63").Value = 3799
$ws.Range("E63").Value = 3777
$ws.Range("F63").Value = 2623
$ws.Range("G63").Value = 1361
$ws.Range("H63").Value = 1176

# Row 64
$ws.Range("B64").Value = 1381
$ws.Range("C64").Value = 1014
$ws.Range("D64").Value = 0.36
$ws.Range("I64").Value = 1381
$ws.Range("J64").Value = 1014

# Row 65
$ws.Range("B65").Value = 1381
$ws.Range("C65").Value = 1014
$ws.Range("D65").Value = 0.36
$ws.Range("I65").Value = 1381
$ws.Range("J65").Value = 1014

# Row 67
$ws.Range("B67").Value = 301417
$ws.Range("C67").Value = 257009
$ws.Range("D67").Value = 0.17
$ws.Range("E67").Value = 27984
$ws.Range("F67").Value = 23092
$ws.Range("G67").Value = 245471
$ws.Range("H67").Value = 206578
$ws.Range("I67").Value = 23173
$ws.Range("J67").Value = 23135
$ws.Range("K67").Value = 4789
$ws.Range("L67").Value = 4204
